# Project Sample Project is saved (rule row B11 label changes from "R40" to "1").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B11 currently holds the text "R40" and must become the text "1" while
# keeping its existing cell formatting (borders/style) untouched.
#
# A direct $ws.Range("B11").Value = "1" (or "'1") would make Excel re-parse
# the literal as user input: a bare "1" is stored as a *number*, and a
# quote-prefixed "'1" is kept as text but forces Excel to stamp the cell
# with a new "number stored as text" style - changing B11's style id.
#
# To land exactly on a *text* "1" without touching B11's style, compute the
# string in a scratch cell via a formula (formula results aren't subject to
# the "looks like a number" input heuristic), copy just the resulting value
# onto B11 (format is left alone), then remove the scratch cell again.
$scratch = $ws.Range("Z1")
$scratch.Formula = "=CHAR(49)"
$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)
$scratch.Clear()
